$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9491214156150818
$ws.Range("B1").Value = 1.647272109985352
$ws.Range("D1").Value = 1.805267333984375
$ws.Range("E1").Value = 1.076490759849548
